$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update header row text -------------------------------------------------
# C1 changes from "CATEGORY TITLE" to "OFFICE/RECIPIENT TITLE"
$ws.Range("C1").Value = "OFFICE/RECIPIENT TITLE"
# H1 gets a new label "Status (usually auto-entered)"
$ws.Range("H1").Value = "Status (usually auto-entered)"

# --- Remove the now-unused column I in the header rows ----------------------
$ws.Range("I1:I8").Clear()

# --- Add a comment on H1 describing how the Status column works ------------
$commentText = @"
Sorkhab, Drake Liu:
This is done automatically by the program.
Examples:
"Previously emailed on 12/22/2023."
"Previous emailed on 12/22/2023. Previously emailed on 12/25/2023."
"Document for 2023 received and awaiting printing."
"Document for 2023 printed and awaiting upload."
"Document for 2023 uploaded."

"@
$comment = $ws.Range("H1").AddComment($commentText)

# --- Update the selected cell in the sheet view -----------------------------
[void]$ws.Activate()
[void]$ws.Range("E2").Select()
